$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 213.66667
$ws.Cells.Item(53, 9).Value = 113.17647
$ws.Cells.Item(53, 10).Value = 457.7143
$ws.Cells.Item(53, 11).Value = 113.17647
$ws.Cells.Item(53, 12).Value = 457.7143
$ws.Cells.Item(53, 13).Value = 523.82353
$ws.Cells.Item(53, 14).Value = -1731.7143
$ws.Cells.Item(86, 8).Value = 50002988
$ws.Cells.Item(86, 9).Value = 2428
$ws.Cells.Item(86, 10).Value = 76926370
$ws.Cells.Item(86, 11).Value = 2428
$ws.Cells.Item(86, 12).Value = 76926370
$ws.Cells.Item(86, 13).Value = -1305
$ws.Cells.Item(86, 14).Value = -76928616
$ws.Cells.Item(89, 8).Value = 50002988
$ws.Cells.Item(89, 9).Value = 2428
$ws.Cells.Item(89, 10).Value = 76926370
$ws.Cells.Item(89, 11).Value = 12140
$ws.Cells.Item(89, 12).Value = 384631850
$ws.Cells.Item(89, 13).Value = -6524
$ws.Cells.Item(89, 14).Value = -384643082
$ws.Cells.Item(98, 8).Value = 1625.3334
$ws.Cells.Item(98, 9).Value = 1260.9565
$ws.Cells.Item(98, 10).Value = 10006
$ws.Cells.Item(98, 11).Value = 1260.9565
$ws.Cells.Item(98, 12).Value = 10006
$ws.Cells.Item(98, 13).Value = 237.0435
$ws.Cells.Item(98, 14).Value = -13002
$ws.Cells.Item(112, 8).Value = 1069.2307
$ws.Cells.Item(112, 9).Value = 1300
$ws.Cells.Item(112, 10).Value = 1000
$ws.Cells.Item(112, 11).Value = 3900
$ws.Cells.Item(112, 12).Value = 3000
$ws.Cells.Item(112, 13).Value = -2792
$ws.Cells.Item(112, 14).Value = -5216
$ws.Cells.Item(122, 8).Value = 1625.3334
$ws.Cells.Item(122, 9).Value = 1260.9565
$ws.Cells.Item(122, 10).Value = 10006
$ws.Cells.Item(122, 11).Value = 3782.8695
$ws.Cells.Item(122, 12).Value = 30018
$ws.Cells.Item(122, 13).Value = -1332.8695
$ws.Cells.Item(122, 14).Value = -34918
$ws.Cells.Item(135, 8).Value = 2618.7827
$ws.Cells.Item(135, 9).Value = 1800.4667
$ws.Cells.Item(135, 11).Value = 16204.2003
$ws.Cells.Item(135, 13).Value = -13669.2003
$ws.Cells.Item(137, 8).Value = 13515024
$ws.Cells.Item(137, 9).Value = 1644.2307
$ws.Cells.Item(137, 10).Value = 45455740
$ws.Cells.Item(137, 11).Value = 4932.6921
$ws.Cells.Item(137, 12).Value = 136367220
$ws.Cells.Item(137, 13).Value = -2382.6921
$ws.Cells.Item(137, 14).Value = -136372320
$ws.Cells.Item(138, 8).Value = 4264.55
$ws.Cells.Item(138, 9).Value = 1338.5294
$ws.Cells.Item(138, 10).Value = 4863.8555
$ws.Cells.Item(138, 11).Value = 4015.5882
$ws.Cells.Item(138, 12).Value = 14591.5665
$ws.Cells.Item(138, 13).Value = 1124.4118
$ws.Cells.Item(138, 14).Value = -24871.5665

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 24566.934
$ws.Cells.Item(32, 9).Value = 20963.928
$ws.Cells.Item(32, 10).Value = 64200
$ws.Cells.Item(32, 11).Value = 20963.928
$ws.Cells.Item(32, 12).Value = 64200
$ws.Cells.Item(32, 13).Value = -20676.928
$ws.Cells.Item(32, 14).Value = -64774

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2797.6
$ws.Cells.Item(86, 9).Value = 2381.3076
$ws.Cells.Item(86, 11).Value = 2381.3076
$ws.Cells.Item(86, 13).Value = -1258.3076
$ws.Cells.Item(89, 8).Value = 2797.6
$ws.Cells.Item(89, 9).Value = 2381.3076
$ws.Cells.Item(89, 11).Value = 11906.538
$ws.Cells.Item(89, 13).Value = -6290.538
$ws.Cells.Item(94, 8).Value = 1825.5
$ws.Cells.Item(94, 9).Value = 1300.7273
$ws.Cells.Item(94, 10).Value = 2980
$ws.Cells.Item(94, 11).Value = 1300.7273
$ws.Cells.Item(94, 12).Value = 2980
$ws.Cells.Item(94, 13).Value = -849.7273
$ws.Cells.Item(94, 14).Value = -3882
$ws.Cells.Item(107, 8).Value = 1614.8
$ws.Cells.Item(107, 9).Value = 1353.6666
$ws.Cells.Item(107, 11).Value = 1353.6666
$ws.Cells.Item(107, 13).Value = 566.3334
$ws.Cells.Item(132, 8).Value = 110000
$ws.Cells.Item(132, 10).Value = 110000
$ws.Cells.Item(132, 12).Value = 110000
$ws.Cells.Item(132, 14).Value = -120120
$ws.Cells.Item(134, 8).Value = 36517.8
$ws.Cells.Item(134, 9).Value = 3210.1904
$ws.Cells.Item(134, 10).Value = 114235.555
$ws.Cells.Item(134, 11).Value = 9630.5712
$ws.Cells.Item(134, 12).Value = 342706.665
$ws.Cells.Item(134, 13).Value = -7095.5712
$ws.Cells.Item(134, 14).Value = -347776.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2228.3428
$ws.Cells.Item(31, 9).Value = 1755.1111
$ws.Cells.Item(31, 10).Value = 2729.4119
$ws.Cells.Item(31, 11).Value = 1755.1111
$ws.Cells.Item(31, 12).Value = 2729.4119
$ws.Cells.Item(31, 13).Value = -1460.1111
$ws.Cells.Item(31, 14).Value = -3319.4119
$ws.Cells.Item(34, 8).Value = 2228.3428
$ws.Cells.Item(34, 9).Value = 1755.1111
$ws.Cells.Item(34, 10).Value = 2729.4119
$ws.Cells.Item(34, 11).Value = 1755.1111
$ws.Cells.Item(34, 12).Value = 2729.4119
$ws.Cells.Item(34, 13).Value = -1553.1111
$ws.Cells.Item(34, 14).Value = -3133.4119
$ws.Cells.Item(140, 8).Value = 54298.89
$ws.Cells.Item(140, 10).Value = 54298.89
$ws.Cells.Item(140, 12).Value = 54298.89
$ws.Cells.Item(140, 14).Value = -64658.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1043.3265
$ws.Cells.Item(68, 10).Value = 1273.6154
$ws.Cells.Item(68, 12).Value = 3820.8462
$ws.Cells.Item(68, 14).Value = -5442.8462
$ws.Cells.Item(71, 8).Value = 1043.3265
$ws.Cells.Item(71, 10).Value = 1273.6154
$ws.Cells.Item(71, 12).Value = 11462.5386
$ws.Cells.Item(71, 14).Value = -19574.5386
$ws.Cells.Item(107, 8).Value = 857.4783
$ws.Cells.Item(107, 10).Value = 1108.5454
$ws.Cells.Item(107, 12).Value = 3325.6362
$ws.Cells.Item(107, 14).Value = -7165.6362
$ws.Cells.Item(113, 8).Value = 217960.7
$ws.Cells.Item(113, 9).Value = 550
$ws.Cells.Item(113, 10).Value = 222792.05
$ws.Cells.Item(113, 11).Value = 1650
$ws.Cells.Item(113, 12).Value = 668376.1499999999
$ws.Cells.Item(113, 13).Value = 520
$ws.Cells.Item(113, 14).Value = -672716.1499999999
$ws.Cells.Item(122, 8).Value = 395.8
$ws.Cells.Item(122, 9).Value = 395.8
$ws.Cells.Item(122, 11).Value = 3562.2
$ws.Cells.Item(122, 13).Value = -1112.2
$ws.Cells.Item(131, 8).Value = 20051.873
$ws.Cells.Item(131, 9).Value = 84649.164
$ws.Cells.Item(131, 10).Value = 2024.721
$ws.Cells.Item(131, 11).Value = 253947.492
$ws.Cells.Item(131, 12).Value = 6074.163
$ws.Cells.Item(131, 13).Value = -248907.492
$ws.Cells.Item(131, 14).Value = -16154.163
$ws.Cells.Item(134, 8).Value = 5357.7827
$ws.Cells.Item(134, 9).Value = 3659.2144
$ws.Cells.Item(134, 11).Value = 10977.6432
$ws.Cells.Item(134, 13).Value = -5907.643199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 9).Value = 4097.615
$ws.Cells.Item(70, 10).Value = 4739.75
$ws.Cells.Item(70, 11).Value = 4097.615
$ws.Cells.Item(70, 12).Value = 4739.75
$ws.Cells.Item(70, 13).Value = -3827.615
$ws.Cells.Item(70, 14).Value = -5279.75
$ws.Cells.Item(73, 9).Value = 4097.615
$ws.Cells.Item(73, 10).Value = 4739.75
$ws.Cells.Item(73, 11).Value = 4097.615
$ws.Cells.Item(73, 12).Value = 4739.75
$ws.Cells.Item(73, 13).Value = -3161.615
$ws.Cells.Item(73, 14).Value = -6611.75
$ws.Cells.Item(138, 8).Value = 19963.625
$ws.Cells.Item(138, 10).Value = 19963.625
$ws.Cells.Item(138, 12).Value = 19963.625
$ws.Cells.Item(138, 14).Value = -30243.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2201.25
$ws.Cells.Item(7, 9).Value = 2200
$ws.Cells.Item(7, 10).Value = 2202.5
$ws.Cells.Item(7, 11).Value = 2200
$ws.Cells.Item(7, 12).Value = 2202.5
$ws.Cells.Item(7, 13).Value = -2088
$ws.Cells.Item(7, 14).Value = -2426.5
$ws.Cells.Item(46, 8).Value = 417
$ws.Cells.Item(46, 9).Value = 450.5
$ws.Cells.Item(46, 10).Value = 350
$ws.Cells.Item(46, 11).Value = 450.5
$ws.Cells.Item(46, 12).Value = 350
$ws.Cells.Item(46, 13).Value = -262.5
$ws.Cells.Item(46, 14).Value = -726
$ws.Cells.Item(93, 8).Value = 1465.7
$ws.Cells.Item(93, 9).Value = 1269.125
$ws.Cells.Item(93, 10).Value = 2252
$ws.Cells.Item(93, 11).Value = 1269.125
$ws.Cells.Item(93, 12).Value = 2252
$ws.Cells.Item(93, 13).Value = -21.125
$ws.Cells.Item(93, 14).Value = -4748
$ws.Cells.Item(103, 8).Value = 12601.5
$ws.Cells.Item(103, 10).Value = 12601.5
$ws.Cells.Item(103, 12).Value = 12601.5
$ws.Cells.Item(103, 14).Value = -14945.5
$ws.Cells.Item(126, 8).Value = 2201.25
$ws.Cells.Item(126, 9).Value = 2200
$ws.Cells.Item(126, 10).Value = 2202.5
$ws.Cells.Item(126, 11).Value = 6600
$ws.Cells.Item(126, 12).Value = 6607.5
$ws.Cells.Item(126, 13).Value = -4130
$ws.Cells.Item(126, 14).Value = -11547.5
$ws.Cells.Item(132, 8).Value = 2528163.2
$ws.Cells.Item(132, 9).Value = 3971084.8
$ws.Cells.Item(132, 10).Value = 3050.25
$ws.Cells.Item(132, 11).Value = 11913254.4
$ws.Cells.Item(132, 12).Value = 9150.75
$ws.Cells.Item(132, 13).Value = -11910724.4
$ws.Cells.Item(132, 14).Value = -14210.75
$ws.Cells.Item(139, 8).Value = 40168.332
$ws.Cells.Item(139, 10).Value = 40168.332
$ws.Cells.Item(139, 12).Value = 40168.332
$ws.Cells.Item(139, 14).Value = -50448.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1198.5
$ws.Cells.Item(100, 9).Value = 2700
$ws.Cells.Item(100, 10).Value = 447.75
$ws.Cells.Item(100, 11).Value = 5400
$ws.Cells.Item(100, 12).Value = 895.5
$ws.Cells.Item(100, 13).Value = -4859
$ws.Cells.Item(100, 14).Value = -1977.5
$ws.Cells.Item(132, 8).Value = 1262.2549
$ws.Cells.Item(132, 9).Value = 791.5333000000001
$ws.Cells.Item(132, 10).Value = 1934.7142
$ws.Cells.Item(132, 11).Value = 2374.5999
$ws.Cells.Item(132, 12).Value = 5804.142599999999
$ws.Cells.Item(132, 13).Value = 155.4000999999998
$ws.Cells.Item(132, 14).Value = -10864.1426
$ws.Cells.Item(138, 8).Value = 40120
$ws.Cells.Item(138, 10).Value = 40120
$ws.Cells.Item(138, 12).Value = 40120
$ws.Cells.Item(138, 14).Value = -50400
